{"js": "// Rewrite the \"Impact\" bullets under \"KEY ACHIEVEMENTS AND IMPACT\" from\n// six job-duty-style bullets down to four impact-focused accomplishment\n// statements, per the commit:\n//   \"Fix Key Achievements to use proper accomplishment statements\"\n//\n// Old bullets (in order):\n//   1. Delivered $4.9M additional revenue through continuous testing and\n//      optimization, increased conversion rates by 23%\n//   2. Built redistricting platform used by thousands of analysts\n//      nationwide with real-time collaborative editing and Census\n//      integration, serving 12,847 analysts across 89 organizations\n//   3. Achieved 87% prediction accuracy for voter turnout vs. industry\n//      standard of 71%, reducing polling error margins from \u00b14.2% to\n//      \u00b12.1%\n//   4. Trigonometric algorithm for boundary estimation reduced mapping\n//      costs by 73.5%, saving campaigns and organizations $4.7M and\n//      enabling smaller nonprofits to conduct analysis\n//   5. Discovered systematic race coding errors affecting all Black and\n//      Asian-American voters, developed geospatial machine learning\n//      algorithms improving classification accuracy from 23% to 64%\n//   6. Developed longitudinal data analysis methods using geospatial\n//      techniques that improved segmentation accuracy by 34% and survey\n//      incidence rates by 28%, reducing polling costs while increasing\n//      response quality\n//\n// New bullets (in order):\n//   1. Platform impact: Built redistricting system serving 12,847\n//      analysts across 89 organizations\n//   2. Real-time collaboration at national scale\n//   3. Revenue generation: Delivered $4.9M additional revenue through\n//      optimization\n//   4. 23% conversion rate improvement\n//\n// Bullets 1-4 get their text swapped in place; bullets 5 and 6 are\n// deleted outright (their paragraphs removed).\n\nconst body = context.document.body;\n\n// Locate the \"KEY ACHIEVEMENTS AND IMPACT\" heading so we scope our search\n// to that section only -- several of these bullet strings (e.g. the\n// \"Achieved 87%...\" and \"Trigonometric algorithm...\" lines) also appear\n// verbatim earlier, under \"PROFESSIONAL EXPERIENCE\", and must be left\n// untouched.\nconst allParagraphs = body.paragraphs;\nallParagraphs.load(\"items/text\");\nawait context.sync();\n\nlet headingParagraph = null;\nfor (let i = 0; i < allParagraphs.items.length; i++) {\n  if (allParagraphs.items[i].text === \"KEY ACHIEVEMENTS AND IMPACT\") {\n    headingParagraph = allParagraphs.items[i];\n    break;\n  }\n}\nif (!headingParagraph) {\n  throw new Error('Could not find \"KEY ACHIEVEMENTS AND IMPACT\" heading');\n}\n\n// Range from just after the heading through the end of the document body;\n// all subsequent search/replace operations are confined to this range.\nconst sectionRange = headingParagraph.getRange(\"After\").expandTo(body.getRange(\"End\"));\n\n// Old bullet text -> new bullet text (null => delete the paragraph).\nconst replacements = [\n  {\n    old: \"\\u2022 Delivered $4.9M additional revenue through continuous testing and optimization, increased conversion rates by 23%\",\n    new: \"\\u2022 Platform impact: Built redistricting system serving 12,847 analysts across 89 organizations\",\n  },\n  {\n    old: \"\\u2022 Built redistricting platform used by thousands of analysts nationwide with real-time collaborative editing and Census integration, serving 12,847 analysts across 89 organizations\",\n    new: \"\\u2022 Real-time collaboration at national scale\",\n  },\n  {\n    old: \"\\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from \\u00b14.2% to \\u00b12.1%\",\n    new: \"\\u2022 Revenue generation: Delivered $4.9M additional revenue through optimization\",\n  },\n  {\n    old: \"\\u2022 Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M and enabling smaller nonprofits to conduct analysis\",\n    new: \"\\u2022 23% conversion rate improvement\",\n  },\n  {\n    old: \"\\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving classification accuracy from 23% to 64%\",\n    new: null,\n  },\n  {\n    old: \"\\u2022 Developed longitudinal data analysis methods using geospatial techniques that improved segmentation accuracy by 34% and survey incidence rates by 28%, reducing polling costs while increasing response quality\",\n    new: null,\n  },\n];\n\n// Resolve each bullet's paragraph up front (search results are stable\n// Range objects even after earlier paragraphs in the list are edited).\nconst searchResults = replacements.map((entry) =>\n  sectionRange.search(entry.old, { matchCase: true, matchWholeWord: false })\n);\nsearchResults.forEach((results) => results.load(\"items\"));\nawait context.sync();\n\nconst targetParagraphs = searchResults.map((results, idx) => {\n  if (results.items.length !== 1) {\n    throw new Error(\n      `Expected exactly 1 match for bullet ${idx + 1}, found ${results.items.length}`\n    );\n  }\n  return results.items[0].paragraphs.getFirst();\n});\nawait context.sync();\n\n// Apply edits: rewrite the first four bullets in place, delete the last two.\nfor (let i = 0; i < replacements.length; i++) {\n  const { new: newText } = replacements[i];\n  const paragraph = targetParagraphs[i];\n  if (newText === null) {\n    paragraph.delete();\n  } else {\n    paragraph.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Rewrite the \"Impact\" bullets under \"KEY ACHIEVEMENTS AND IMPACT\" from six\n# job-duty-style bullets down to four impact-focused accomplishment\n# statements, per the commit:\n#   \"Fix Key Achievements to use proper accomplishment statements\"\n#\n# Old bullets (in order):\n#   1. Delivered $4.9M additional revenue through continuous testing and\n#      optimization, increased conversion rates by 23%\n#   2. Built redistricting platform used by thousands of analysts\n#      nationwide with real-time collaborative editing and Census\n#      integration, serving 12,847 analysts across 89 organizations\n#   3. Achieved 87% prediction accuracy for voter turnout vs. industry\n#      standard of 71%, reducing polling error margins from \u00b14.2% to\n#      \u00b12.1%\n#   4. Trigonometric algorithm for boundary estimation reduced mapping\n#      costs by 73.5%, saving campaigns and organizations $4.7M and\n#      enabling smaller nonprofits to conduct analysis\n#   5. Discovered systematic race coding errors affecting all Black and\n#      Asian-American voters, developed geospatial machine learning\n#      algorithms improving classification accuracy from 23% to 64%\n#   6. Developed longitudinal data analysis methods using geospatial\n#      techniques that improved segmentation accuracy by 34% and survey\n#      incidence rates by 28%, reducing polling costs while increasing\n#      response quality\n#\n# New bullets (in order):\n#   1. Platform impact: Built redistricting system serving 12,847\n#      analysts across 89 organizations\n#   2. Real-time collaboration at national scale\n#   3. Revenue generation: Delivered $4.9M additional revenue through\n#      optimization\n#   4. 23% conversion rate improvement\n#\n# Bullets 1-4 get their text swapped in place; bullets 5 and 6 are deleted\n# outright (their paragraphs removed).\n\n$d = $word.ActiveDocument\n\n$bullet = [char]0x2022\n$plusMinus = [char]0xB1\n\n# Locate the \"KEY ACHIEVEMENTS AND IMPACT\" heading so every search below is\n# scoped to that section only -- several of these bullet strings (e.g. the\n# \"Achieved 87%...\" and \"Trigonometric algorithm...\" lines) also appear\n# verbatim earlier, under \"PROFESSIONAL EXPERIENCE\", and must be left\n# untouched.\n$headingIdx = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13) -eq \"KEY ACHIEVEMENTS AND IMPACT\") {\n        $headingIdx = $i\n        break\n    }\n}\nif ($headingIdx -eq -1) {\n    throw \"Could not find the 'KEY ACHIEVEMENTS AND IMPACT' heading\"\n}\n\n# Old bullet text -> new bullet text ($null => delete the paragraph).\n$replacements = @(\n    @{\n        Old = \"$bullet Delivered `$4.9M additional revenue through continuous testing and optimization, increased conversion rates by 23%\"\n        New = \"$bullet Platform impact: Built redistricting system serving 12,847 analysts across 89 organizations\"\n    },\n    @{\n        Old = \"$bullet Built redistricting platform used by thousands of analysts nationwide with real-time collaborative editing and Census integration, serving 12,847 analysts across 89 organizations\"\n        New = \"$bullet Real-time collaboration at national scale\"\n    },\n    @{\n        Old = \"$bullet Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from ${plusMinus}4.2% to ${plusMinus}2.1%\"\n        New = \"$bullet Revenue generation: Delivered `$4.9M additional revenue through optimization\"\n    },\n    @{\n        Old = \"$bullet Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M and enabling smaller nonprofits to conduct analysis\"\n        New = \"$bullet 23% conversion rate improvement\"\n    },\n    @{\n        Old = \"$bullet Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving classification accuracy from 23% to 64%\"\n        New = $null\n    },\n    @{\n        Old = \"$bullet Developed longitudinal data analysis methods using geospatial techniques that improved segmentation accuracy by 34% and survey incidence rates by 28%, reducing polling costs while increasing response quality\"\n        New = $null\n    }\n)\n\nforeach ($r in $replacements) {\n    # Re-derive the section range on every iteration: earlier deletions /\n    # replacements in this loop shift character offsets, so the heading's\n    # paragraph index (which does not move) is the only thing we trust.\n    $heading = $d.Paragraphs.Item($headingIdx)\n    $sectionRange = $d.Range($heading.Range.End, $d.Content.End)\n\n    $found = $sectionRange.Find.Execute($r.Old)\n    if (-not $found) {\n        throw \"Could not find expected bullet text: $($r.Old)\"\n    }\n\n    if ($null -eq $r.New) {\n        # Delete the whole paragraph (including its end-of-paragraph mark)\n        # via the document's Paragraphs collection -- deleting the found\n        # Range alone leaves a blank paragraph behind because Find's match\n        # stops short of the paragraph mark.\n        $paraIndex = $sectionRange.Paragraphs.First.Index\n        $d.Paragraphs.Item($paraIndex).Range.Delete()\n    } else {\n        $sectionRange.Text = $r.New\n    }\n}\n"}
